$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.978.09'
$ws.Range('E2').Value = '  +1.44%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.180.63'
$ws.Range('E3').Value = '  +0.69%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '594.91'
$ws.Range('E5').Value = '  +3.60%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '153.36'
$ws.Range('E6').Value = '  +2.18%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.178.89'
$ws.Range('E8').Value = '  +0.69%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.542'
$ws.Range('E9').Value = '  +2.65%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.159'
$ws.Range('E10').Value = '  -0.82%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.08'
$ws.Range('E11').Value = '  -0.10%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.515'
$ws.Range('E12').Value = '  +3.60%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000267'
$ws.Range('E13').Value = '  +0.49%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '38.91'
$ws.Range('E14').Value = '  +4.66%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.699.19'
$ws.Range('E15').Value = '  +0.61%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.981.09'
$ws.Range('E16').Value = '  +1.38%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.43'
$ws.Range('E17').Value = '  +4.67%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.176.66'
$ws.Range('E18').Value = '  +0.62%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.112'
$ws.Range('E19').Value = '  +0.65%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '509.31'
$ws.Range('E20').Value = '  +0.67%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.34'
$ws.Range('E21').Value = '  +3.40%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.737'
$ws.Range('E22').Value = '  +2.45%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.02'
$ws.Range('E23').Value = '  +3.73%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '15.02'
$ws.Range('E24').Value = '  -1.91%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.84'
$ws.Range('E25').Value = '  +0.57%  '
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.27'
$ws.Range('E27').Value = '  +4.07%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.98'
$ws.Range('E28').Value = '  +2.53%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.29'
$ws.Range('E29').Value = '  +5.72%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.99'
$ws.Range('E30').Value = '  +13.32%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.87'
$ws.Range('E31').Value = '  +3.24%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '28.11'
$ws.Range('E32').Value = '  +1.82%  '
$ws.Range('E33').Value = '  +2.50%  '
$ws.Range('E34').Value = '  +0.06%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.48'
$ws.Range('E35').Value = '  -0.60%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '54.79'
$ws.Range('E36').Value = '  -0.23%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0902'
$ws.Range('E37').Value = '  +0.41%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '481.12'
$ws.Range('E38').Value = '  +3.11%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0419'
$ws.Range('E39').Value = '  -0.22%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.82'
$ws.Range('E40').Value = '  +1.68%  '
$ws.Range('B41').Value = 'TheGraph'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.301'
$ws.Range('E41').Value = '  +6.32%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.122'
$ws.Range('E42').Value = '  +3.65%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.84'
$ws.Range('E43').Value = '  -4.73%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0₃0656'
$ws.Range('E44').Value = '  +11.86%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.899.98'
$ws.Range('E45').Value = '  -4.81%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.42'
$ws.Range('E46').Value = '  -1.01%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '28.42'
$ws.Range('E47').Value = '  -0.36%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.117'
$ws.Range('E49').Value = '  +2.01%  '
$ws.Range('E50').Value = '  +2.58%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.59'
$ws.Range('E51').Value = '  +5.49%  '
